# The post that used to live on row 154 ("「理想的なスイカの選び方」...") was
# removed from the source data. Deleting its entire row shifts every
# subsequent row (155-205) up by one, which is exactly what the target
# diff shows (old row 155 becomes new row 154, ..., old row 205 becomes
# new row 204) and shrinks the used range from A1:C205 to A1:C204.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(154).Delete()
